$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1078
$ws.Range("F3").Value = 4627
$ws.Range("F5").Value = 171
$ws.Range("F6").Value = 1799
$ws.Range("F7").Value = 39
$ws.Range("F8").Value = 725
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 407
$ws.Range("F13").Value = 1568
$ws.Range("F14").Value = 806
$ws.Range("F15").Value = 766
$ws.Range("F19").Value = 160
$ws.Range("F22").Value = 387
$ws.Range("F23").Value = 2502
$ws.Range("F25").Value = 1538
$ws.Range("F26").Value = 486
$ws.Range("F27").Value = 526
$ws.Range("F29").Value = 4218

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 360
$ws.Range("F12").Value = 34
$ws.Range("F18").Value = 280
$ws.Range("F20").Value = 140

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1326
$ws.Range("F5").Value = 1724
$ws.Range("F6").Value = 1070
$ws.Range("F7").Value = 235

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1326
$ws.Range("F3").Value = 1724
$ws.Range("F4").Value = 1070
$ws.Range("F5").Value = 235
$ws.Range("F6").Value = 1078
$ws.Range("F8").Value = 4627
$ws.Range("F10").Value = 171
$ws.Range("F11").Value = 1799
$ws.Range("F12").Value = 725
$ws.Range("F13").Value = 360
$ws.Range("F15").Value = 13
$ws.Range("F16").Value = 407
$ws.Range("F18").Value = 1568
$ws.Range("F21").Value = 806
$ws.Range("F22").Value = 766
$ws.Range("F26").Value = 160
$ws.Range("F29").Value = 280
$ws.Range("F34").Value = 387
$ws.Range("F35").Value = 140
$ws.Range("F37").Value = 2502
$ws.Range("F43").Value = 1538
$ws.Range("F44").Value = 486
$ws.Range("F45").Value = 526
$ws.Range("F48").Value = 4218
